# Applies the "Quantum" -> "Mathematics" themed rewrite described by the
# commit diff, using the Word COM/object model only.

function Set-ParaText {
    param($doc, $index, $text)
    $p = $doc.Paragraphs.Item($index)
    $r = $p.Range
    $r2 = $doc.Range($r.Start, $r.End - 1)
    $r2.Text = $text
}

$d = $word.ActiveDocument

# --- 1. Title ---
Set-ParaText $d 1 "A Journey Through the Magic of Mathematics"

# --- 2. Author name ---
Set-ParaText $d 2 "Samuel Davies"

# --- 3. Author email ---
Set-ParaText $d 3 "samueldavies@edumail.org"

# --- 4. Paragraph 4 is the blank spacer paragraph - left untouched ---

# --- 5. Main body paragraph (paragraph 5): three sentences separated by
#        two manual line breaks (chr 11 == vertical tab == <w:br/>) ---
$br = [char]11

$bodyPart1 = "Mathematics, the universal language of the universe, beckons us on an awe-inspiring intellectual adventure." +
    " This enchanting realm has captured the imagination of brilliant minds throughout history, inspiring discoveries that have shaped our understanding of the world." +
    " From ancient civilizations to modern-day frontiers, mathematics weaves intricate threads connecting science, technology, engineering, and art." +
    " It is the orchestra conductor of our universe, orchestrating the rhythm and melody of existence."

$bodyPart2 = "In the symphony of mathematics, we witness the harmony of patterns, the elegance of symmetry, and the power of logical reasoning." +
    " It holds the key to unlocking nature's deepest secrets, empowering us to unveil the enigmas of the cosmos." +
    " Mathematics empowers us to navigate the complexity of modern life, from financial transactions to intricate engineering marvels." +
    " It is the language of innovation, unraveling new frontiers of knowledge and shaping the trajectory of human progress."

$bodyPart3 = "Mathematics is more than a subject; it's a way of thinking, a lens through which we can perceive the world with greater clarity." +
    " It cultivates analytical and problem-solving skills, nurturing critical thinking and inspiring creativity." +
    " By delving into the depths of mathematics, we cultivate a mindset that embraces precision, logic, and imagination, enabling us to become more effective problem solvers and informed decision-makers."

$bodyFull = $bodyPart1 + $br + $br + $bodyPart2 + $br + $br + $bodyPart3

Set-ParaText $d 5 $bodyFull

# --- 6. Paragraph 6 "Summary" heading - left untouched ---

# --- 7. Summary paragraph (paragraph 7) ---
$summaryFull = "This exploration of mathematics highlights its remarkable allure, unveiling its role as the universal language underpinning the fabric of our universe." +
    " Mathematics weaves enchanting connections between the worlds of science, technology, and art, inspiring discoveries that have profoundly shaped our understanding of existence." +
    " It empowers us to unlock nature's secrets, navigate the complexities of modern life, and cultivate analytical minds capable of solving real-world problems." +
    " Embracing the enchantment of " +
    "mathematics, we transform into informed thinkers and effective problem-solvers, ready to contribute to a future shaped by innovation and discovery."

Set-ParaText $d 7 $summaryFull

# --- 8. A new, empty paragraph is appended after the summary paragraph ---
$d.Paragraphs.Add() | Out-Null

Write-Host "edit complete"
